$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 94 was an empty separator row; deleting it shifts rows 95:106 up to
# 94:105 (the old duplicate row 106 content lands on 105, and the sheet
# shrinks from A1:I106 to A1:I105).
$ws.Rows(94).Delete()

# The hidden _xlnm._FilterDatabase name tracked a stale, too-small range;
# grow it to cover the full table again.
$names = $wb.Names
$fdb = $names.Item("ECB_Distinta SEF Articoli !_FilterDatabase")
$fdb.RefersTo = "='ECB_Distinta SEF Articoli '!`$A`$1:`$I`$105"

# Grow the visible AutoFilter range to match the full table as well.
$ws.AutoFilterMode = $false
$ws.Range("A1:I105").AutoFilter()

# Restore the user's on-sheet selection recorded after the edit.
$ws.Range("D108").Select()
